# This workbook contains a single data table (rows 2..257) describing daily
# price records for "Ciboulette" at "Feria Lagunitas de Puerto Montt".
# The edit adds one new record at the top of the series (row 127, which is
# where new rows get inserted chronologically in this particular sheet),
# pushing every subsequent record (rows 127..257) down by one row, and a new
# row 258 is appended to hold what used to be the last record (old row 257).
#
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# O (Origen) and P (Precio $/Kg) all shift down by one row for rows
# 128..258. Row 127 keeps its former K..P values but gets a brand new
# Fecha (D) and Volumen (J). Columns A, B, C, E, F, G, H, I, Q, R are
# constant for every data row in this sheet, so shifting them along with
# the rest of the block is harmless, and the new row 258 just needs them
# copied over from row 257.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the whole block (D..P) for the rows that will move down by one.
$src = $ws.Range("D127:P257").Value()

# Make sure the brand-new row 258 uses the same date format as the rest of
# the Fecha column before we put a date value into it, so Excel doesn't
# invent a new/duplicate number format style for it.
$ws.Range("D258").NumberFormat = $ws.Range("D257").NumberFormat()

# Shift rows 127..257 down to 128..258.
$ws.Range("D128:P258").Value = $src

# Row 127 becomes the new record: new Fecha and Volumen values, while
# K..P (price/unit/origin columns) stay exactly as they were.
$ws.Range("D127").Value = 44810
$ws.Range("J127").Value = 240

# Row 258 also needs the columns that are outside the D:P block and are
# constant across the whole table; just copy them from row 257.
$ws.Range("A258").Value = $ws.Range("A257").Value()
$ws.Range("B258").Value = $ws.Range("B257").Value()
$ws.Range("C258").Value = $ws.Range("C257").Value()
$ws.Range("Q258").Value = $ws.Range("Q257").Value()
$ws.Range("R258").Value = $ws.Range("R257").Value()
